$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the coin table (columns B-E) to the refreshed market snapshot.
# For numeric-looking Price strings (e.g. "1.00", "512.76") we briefly force
# a text NumberFormat so Excel does not auto-convert them to actual numbers,
# then clear the temporary format so the cell keeps the sheet default style
# (matching the original plain-text cells, which carry no style index).

# Row 2
$ws.Range("D2").Value = '74.923.65'
$ws.Range("E2").Value = '  +1.43%  '

# Row 3
$ws.Range("D3").Value = '2.823.29'
$ws.Range("E3").Value = '  +7.80%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '188.33'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.13%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '594.89'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.45%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("E8").Value = '  +3.41%  '

# Row 9
$ws.Range("E9").Value = '  -3.27%  '

# Row 10
$ws.Range("D10").Value = '2.820.73'
$ws.Range("E10").Value = '  +7.81%  '

# Row 11
$ws.Range("E11").Value = '  -1.05%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.371'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.49%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.85'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.99%  '

# Row 14
$ws.Range("D14").Value = '3.346.08'
$ws.Range("E14").Value = '  +8.76%  '

# Row 15
$ws.Range("D15").Value = '74.941.88'
$ws.Range("E15").Value = '  +1.61%  '

# Row 16
$ws.Range("E16").Value = '  -0.42%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.93'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.32%  '

# Row 18
$ws.Range("D18").Value = '2.826.07'
$ws.Range("E18").Value = '  +8.02%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.91'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.69%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.33'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.15%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.20'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.17%  '

# Row 22
$ws.Range("E22").Value = '  -0.67%  '

# Row 23
$ws.Range("E23").Value = '  -0.28%  '

# Row 24
$ws.Range("E24").Value = '  -0.18%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.79'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.59%  '

# Row 26
$ws.Range("D26").Value = '2.953.81'
$ws.Range("E26").Value = '  +8.01%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.16'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.80%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.65'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.65%  '

# Row 29
$ws.Range("E29").Value = '  +11.62%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.62%  '

# Row 31
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '512.76'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.52%  '

# Row 32
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.39'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.12%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.74'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.90%  '

# Row 34
$ws.Range("E34").Value = '  +3.32%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.04%  '

# Row 36
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.12'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.54%  '

# Row 37
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.00'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.38%  '

# Row 38
$ws.Range("E38").Value = '  -0.87%  '

# Row 39
$ws.Range("E39").Value = '  +0.74%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '185.31'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +14.61%  '

# Row 41
$ws.Range("E41").Value = '  +0.02%  '

# Row 42
$ws.Range("E42").Value = '  +4.73%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.02'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.04%  '

# Row 44
$ws.Range("E44").Value = '  +0.31%  '

# Row 45
$ws.Range("E45").Value = '  +3.06%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.01'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.66%  '

# Row 47
$ws.Range("E47").Value = '  +0.35%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0856'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.84%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.573'
$ws.Range("D49").ClearFormats()

# Row 50
$ws.Range("E50").Value = '  +2.44%  '

# Row 51
$ws.Range("E51").Value = '  +8.45%  '
